$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in B1 ("Chi tiết sản phẩm" -> "Chi tiết mặt hàng")
$ws.Range("B1").Value = "Chi tiết mặt hàng"

# Update the column label in B2 ("MaSP" -> "MaMH")
$ws.Range("B2").Value = "MaMH"

# Move the active selection from D12 to B1
$ws.Range("B1").Select()
